$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the extra (4th) row so the sheet shrinks to 3 data rows ---
$ws.Rows(4).Delete()

# --- Row 1: header labels ---
$ws.Range("A1").Value = "varUsername"
$ws.Range("B1").Value = "varPassword"
$ws.Range("C1").Value = "varMulai"
$ws.Range("D1").Value = "varSelesai"
$ws.Range("E1").Value = "varSearch"
$ws.Range("F1").Value = "varUbahSelesai"
$ws.Range("G1").Value = "varBenar"

# --- Row 2 ---
$ws.Range("A2").Value = "gilank.rangesti"
$ws.Range("B2").Value = "Password2"
$ws.Range("C2").Value = 202009
$ws.Range("D2").Value = 202109
$ws.Range("E2").Value = 202009
$ws.Range("F2").Value = 202110
$ws.Range("G2").Value = "Y"

# --- Row 3 ---
$ws.Range("A3").Value = "gilank.rangesti"
$ws.Range("B3").Value = "Password2"
$ws.Range("C3").Value = 202008
$ws.Range("D3").Value = 202108
$ws.Range("E3").Value = "HJKL"
$ws.Range("F3").Value = 202109
$ws.Range("G3").Value = "N"

# --- Match the existing centred cell style on the whole used range ---
$ws.Range("A1:G3").HorizontalAlignment = -4108

# --- Column widths for the newly used columns C:F (closest reachable values) ---
$ws.Range("C1").EntireColumn.ColumnWidth = 14 - (5.0/6.0)
$ws.Range("D1").EntireColumn.ColumnWidth = 15.85546875 - (5.0/6.0)
$ws.Range("E1").EntireColumn.ColumnWidth = 14.7109375 - (5.0/6.0)
$ws.Range("F1").EntireColumn.ColumnWidth = 22.28515625 - (5.0/6.0)

# --- Leave the cursor on C5, matching the saved selection ---
$ws.Range("C5").Select()
